# Edit script: rename sheet, update headers, add two new columns (E,F),
# and rewrite the data values for B:F across all data rows (2-41).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet from "Quarterly" to "Income_Statement"
$ws.Name = "Income_Statement"

# 2. Update header row (row 1): keep A1/B1, rename C1/D1, add E1/F1
$ws.Range("C1").Value = "Cost_of_Goods_Sold"
$ws.Range("D1").Value = "Operating_Expenses"
$ws.Range("E1").Value = "Operating_Profit"
$ws.Range("F1").Value = "Net_Profit"

# Apply the same header formatting (bold, centered, bordered) used by the
# existing header cells to the two newly-added header cells.
$ws.Range("E1:F1").Font.Bold = $true
$ws.Range("E1:F1").HorizontalAlignment = -4108
$ws.Range("E1:F1").VerticalAlignment = -4160
$ws.Range("E1:F1").Borders.LineStyle = 1

$newData = @(
    @(2, 504967.14, 325924.83, 107032.54, 72009.77, 57607.82),
    @(3, 503544.06, 310319.06, 136623.83, 56601.17, 45280.94),
    @(4, 516378.84, 330960.02, 142561.1, 42857.72, 34286.17),
    @(5, 581649.14, 368895.91, 148976.47, 63776.77, 51021.41),
    @(6, 517658.47, 305255.43, 143441.46, 68961.57000000001, 55169.26),
    @(7, 522782.4, 333682.14, 130371.25, 58729.01, 46983.21),
    @(8, 546090.16, 280293.81, 137763.96, 128032.39, 102425.91),
    @(9, 596749.95, 315917.89, 144863.5, 135968.56, 108774.85),
    @(10, 536105.26, 271689.62, 108583.78, 155831.86, 124665.48),
    @(11, 551554.3199999999, 302692.76, 116261.66, 132599.9, 106079.92),
    @(12, 546875.77, 305321.62, 111093.94, 130460.22, 104368.17),
    @(13, 607981.33, 328736.93, 160288.83, 118955.56, 95164.45),
    @(14, 564851.62, 352642.87, 130726.77, 81481.98, 65185.58),
    @(15, 548841.0699999999, 303790.67, 137680.66, 107369.74, 85895.78999999999),
    @(16, 556321.17, 301604.06, 161754.08, 92963.03, 74370.42999999999),
    @(17, 631521.29, 367169.27, 142047.59, 122304.43, 97843.55),
    @(18, 574800.97, 299550.99, 138549.04, 136700.93, 109360.75),
    @(19, 593835.3, 368373.58, 163634.35, 61827.36, 49461.89),
    @(20, 587432.92, 300285.49, 130926.94, 156220.49, 124976.39),
    @(21, 648506.9, 420253.9, 134693.58, 93559.42, 74847.53999999999),
    @(22, 622982.9399999999, 383655.77, 142647.61, 96679.56, 77343.64999999999),
    @(23, 612062.78, 324275.36, 132280.31, 155507.11, 124405.69),
    @(24, 621048.97, 311038.91, 181948.57, 128061.49, 102449.19),
    @(25, 674887.65, 419995.55, 189516.58, 65375.52, 52300.42),
    @(26, 627215.6800000001, 380110.64, 165171.21, 81933.83, 65547.06),
    @(27, 640002.59, 389986.26, 183774.22, 66242.10000000001, 52993.68),
    @(28, 633678.7, 390149.99, 177662.73, 65865.98, 52692.79),
    @(29, 720457.52, 368230.66, 157533.08, 194693.77, 155755.02),
    @(30, 651959.5, 361035.52, 188583.13, 102340.85, 81872.67999999999),
    @(31, 661532.16, 342263.75, 167985.66, 151282.76, 121026.21),
    @(32, 664979.12, 418581.42, 186688.91, 59708.79, 47767.03),
    @(33, 763891.34, 453365.48, 221229.91, 89295.96000000001, 71436.77),
    @(34, 684149.55, 376032.34, 158586.1, 149531.11, 119624.89),
    @(35, 680449.95, 346712.22, 143578.47, 190159.26, 152127.41),
    @(36, 706061.48, 385966.63, 157305.92, 162788.92, 130231.14),
    @(37, 762974.86, 418703.44, 185182.22, 159089.2, 127271.36),
    @(38, 713744.54, 434985.14, 201134.27, 77625.14, 62100.11),
    @(39, 699071.4399999999, 416390.45, 199985.5, 82695.48, 66156.39),
    @(40, 712467.61, 451050.36, 142988.84, 118428.41, 94742.73),
    @(41, 808159.25, 461323.35, 202908.36, 143927.53, 115142.02),
)

# 3. Update the data rows: write the new B/C/D values and the new E/F values
foreach ($row in $newData) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]  # B: Revenue
    $ws.Cells.Item($r, 3).Value = $row[2]  # C: Cost_of_Goods_Sold
    $ws.Cells.Item($r, 4).Value = $row[3]  # D: Operating_Expenses
    $ws.Cells.Item($r, 5).Value = $row[4]  # E: Operating_Profit
    $ws.Cells.Item($r, 6).Value = $row[5]  # F: Net_Profit
}

Write-Output "Edit applied successfully"
